$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2030")

# Update sheet view: zoom level, scroll position, and selection to match target state
$ws.Activate()
$excel.ActiveWindow.Zoom = 145
$ws.Range("A73").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L85").Select() | Out-Null

# Rows where the player order within a position group changed (full row overwrite)
# Row 12: Keke Kouandjio
$ws.Range("A12").Value = "Keke"
$ws.Range("B12").Value = "Kouandjio"
$ws.Range("C12").Value = "SO"
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = "WR"
$ws.Range("F12").Value = "WR"
$ws.Range("G12").Value = "SP"
$ws.Range("H12").Value = "OFF"
$ws.Range("I12").Value = "deep_threat"
$ws.Range("J12").Value = "impact"
$ws.Range("K12").Value = 79
$ws.Range("L12").Value = 81

# Row 13: Demarco Leon
$ws.Range("A13").Value = "Demarco"
$ws.Range("B13").Value = "Leon"
$ws.Range("C13").Value = "JR"
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = "WR"
$ws.Range("F13").Value = "WR"
$ws.Range("G13").Value = "SP"
$ws.Range("H13").Value = "OFF"
$ws.Range("I13").Value = "physical"
$ws.Range("J13").Value = "impact"
$ws.Range("K13").Value = 80
$ws.Range("L13").Value = 80

# Row 14: Finau Lutui
$ws.Range("A14").Value = "Finau"
$ws.Range("B14").Value = "Lutui"
$ws.Range("C14").Value = "JR"
$ws.Range("D14").Value = $true
$ws.Range("E14").Value = "WR"
$ws.Range("F14").Value = "WR"
$ws.Range("G14").Value = "SP"
$ws.Range("H14").Value = "OFF"
$ws.Range("I14").Value = "route_runner"
$ws.Range("J14").Value = "star"
$ws.Range("K14").Value = 79
$ws.Range("L14").Value = 79

# Row 60: Gregory Fields
$ws.Range("A60").Value = "Gregory"
$ws.Range("B60").Value = "Fields"
$ws.Range("C60").Value = "FR"
$ws.Range("D60").Value = $false
$ws.Range("E60").Value = "LOLB"
$ws.Range("F60").Value = "OLB"
$ws.Range("G60").Value = "LB"
$ws.Range("H60").Value = "DEF"
$ws.Range("I60").Value = "power_rusher"
$ws.Range("J60").Value = "elite"
$ws.Range("K60").Value = 76
$ws.Range("L60").Value = 78

# Row 61: Dontrell East
$ws.Range("A61").Value = "Dontrell"
$ws.Range("B61").Value = "East"
$ws.Range("C61").Value = "FR"
$ws.Range("D61").Value = $false
$ws.Range("E61").Value = "LOLB"
$ws.Range("F61").Value = "OLB"
$ws.Range("G61").Value = "LB"
$ws.Range("H61").Value = "DEF"
$ws.Range("I61").Value = "run_stopper"
$ws.Range("J61").Value = "star"
$ws.Range("K61").Value = 75
$ws.Range("L61").Value = 75

# Row 69: Daquan Draper
$ws.Range("A69").Value = "Daquan"
$ws.Range("B69").Value = "Draper"
$ws.Range("C69").Value = "JR"
$ws.Range("D69").Value = $true
$ws.Range("E69").Value = "CB"
$ws.Range("F69").Value = "CB"
$ws.Range("G69").Value = "DB"
$ws.Range("H69").Value = "DEF"
$ws.Range("I69").Value = "slot"
$ws.Range("J69").Value = "star"
$ws.Range("K69").Value = 87
$ws.Range("L69").Value = 88

# Row 70: Cris Small
$ws.Range("A70").Value = "Cris"
$ws.Range("B70").Value = "Small"
$ws.Range("C70").Value = "SO"
$ws.Range("D70").Value = $true
$ws.Range("E70").Value = "CB"
$ws.Range("F70").Value = "CB"
$ws.Range("G70").Value = "DB"
$ws.Range("H70").Value = "DEF"
$ws.Range("I70").Value = "zone"
$ws.Range("J70").Value = "star"
$ws.Range("K70").Value = 87
$ws.Range("L70").Value = 87

# Remaining rows: only fill in column L (overall_end)
$ws.Range("L2").Value = 87
$ws.Range("L3").Value = 82
$ws.Range("L4").Value = 77
$ws.Range("L5").Value = 75
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 87
$ws.Range("L8").Value = 85
$ws.Range("L9").Value = 78
$ws.Range("L10").Value = 73
$ws.Range("L11").Value = 84
$ws.Range("L15").Value = 78
$ws.Range("L16").Value = 76
$ws.Range("L17").Value = 68
$ws.Range("L18").Value = 68
$ws.Range("L19").Value = 66
$ws.Range("L20").Value = 83
$ws.Range("L21").Value = 74
$ws.Range("L22").Value = 71
$ws.Range("L23").Value = 67
$ws.Range("L24").Value = 66
$ws.Range("L25").Value = 65
$ws.Range("L26").Value = 84
$ws.Range("L27").Value = 83
$ws.Range("L28").Value = 75
$ws.Range("L29").Value = 75
$ws.Range("L30").Value = 93
$ws.Range("L31").Value = 88
$ws.Range("L32").Value = 79
$ws.Range("L33").Value = 68
$ws.Range("L34").Value = 65
$ws.Range("L35").Value = 88
$ws.Range("L36").Value = 81
$ws.Range("L37").Value = 75
$ws.Range("L38").Value = 90
$ws.Range("L39").Value = 80
$ws.Range("L40").Value = 74
$ws.Range("L41").Value = 68
$ws.Range("L42").Value = 64
$ws.Range("L43").Value = 88
$ws.Range("L44").Value = 86
$ws.Range("L45").Value = 92
$ws.Range("L46").Value = 83
$ws.Range("L47").Value = 81
$ws.Range("L48").Value = 78
$ws.Range("L49").Value = 89
$ws.Range("L50").Value = 86
$ws.Range("L51").Value = 80
$ws.Range("L52").Value = 73
$ws.Range("L53").Value = 73
$ws.Range("L54").Value = 90
$ws.Range("L55").Value = 89
$ws.Range("L56").Value = 83
$ws.Range("L57").Value = 76
$ws.Range("L58").Value = 71
$ws.Range("L59").Value = 90
$ws.Range("L62").Value = 86
$ws.Range("L63").Value = 85
$ws.Range("L64").Value = 81
$ws.Range("L65").Value = 88
$ws.Range("L66").Value = 80
$ws.Range("L67").Value = 76
$ws.Range("L68").Value = 89
$ws.Range("L71").Value = 83
$ws.Range("L72").Value = 83
$ws.Range("L73").Value = 83
$ws.Range("L74").Value = 77
$ws.Range("L75").Value = 88
$ws.Range("L76").Value = 83
$ws.Range("L77").Value = 82
$ws.Range("L78").Value = 75
$ws.Range("L79").Value = 67
$ws.Range("L80").Value = 80
$ws.Range("L81").Value = 74
$ws.Range("L82").Value = 75
$ws.Range("L83").Value = 86
$ws.Range("L84").Value = 67
$ws.Range("L85").Value = 78
$ws.Range("L86").Value = 66
